# Auto-generated edit script: update cached price/profit values in each sheet's table
# (values refreshed by the scheduled market-data runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4856.2
$ws.Range("I132").Value = 4889.294
$ws.Range("J132").Value = 4668.6665
$ws.Range("K132").Value = 14667.882
$ws.Range("L132").Value = 14005.9995
$ws.Range("M132").Value = -12137.882
$ws.Range("N132").Value = -19065.9995
$ws.Range("H134").Value = 46994.5
$ws.Range("J134").Value = 46994.5
$ws.Range("L134").Value = 46994.5
$ws.Range("N134").Value = -57134.5
$ws.Range("H135").Value = 16668623
$ws.Range("I135").Value = 546.12
$ws.Range("K135").Value = 4915.08
$ws.Range("M135").Value = -2380.08
$ws.Range("H138").Value = 47622060
$ws.Range("I138").Value = 71430390
$ws.Range("J138").Value = 5400
$ws.Range("K138").Value = 214291170
$ws.Range("L138").Value = 16200
$ws.Range("M138").Value = -214286030
$ws.Range("N138").Value = -26480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4627.1567
$ws.Range("I32").Value = 4297.553
$ws.Range("K32").Value = 4297.553
$ws.Range("M32").Value = -4010.553
$ws.Range("H61").Value = 4466.5293
$ws.Range("I61").Value = 4752.643
$ws.Range("K61").Value = 4752.643
$ws.Range("M61").Value = -4540.643
$ws.Range("H63").Value = 3126175.5
$ws.Range("I63").Value = 1306.1111
$ws.Range("K63").Value = 1306.1111
$ws.Range("M63").Value = -620.1111000000001
$ws.Range("H66").Value = 3126175.5
$ws.Range("I66").Value = 1306.1111
$ws.Range("K66").Value = 6530.5555
$ws.Range("M66").Value = -3098.5555
$ws.Range("H74").Value = 50002460
$ws.Range("I74").Value = 125002030
$ws.Range("K74").Value = 125002030
$ws.Range("M74").Value = -125001156
$ws.Range("H77").Value = 50002460
$ws.Range("I77").Value = 125002030
$ws.Range("K77").Value = 625010150
$ws.Range("M77").Value = -625005782
$ws.Range("H97").Value = 200002240
$ws.Range("I97").Value = 2795
$ws.Range("J97").Value = 1000000000
$ws.Range("K97").Value = 2795
$ws.Range("L97").Value = 1000000000
$ws.Range("M97").Value = -2299
$ws.Range("N97").Value = -1000000992
$ws.Range("H132").Value = 20256.75
$ws.Range("I132").Value = 2151.3333
$ws.Range("J132").Value = 74573
$ws.Range("K132").Value = 6453.999899999999
$ws.Range("L132").Value = 223719
$ws.Range("M132").Value = -3923.999899999999
$ws.Range("N132").Value = -228779
$ws.Range("H136").Value = 4466.5293
$ws.Range("I136").Value = 4752.643
$ws.Range("K136").Value = 14257.929
$ws.Range("M136").Value = -11707.929

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4169350.8
$ws.Range("I105").Value = 2535
$ws.Range("J105").Value = 8336166.5
$ws.Range("K105").Value = 2535
$ws.Range("L105").Value = 8336166.5
$ws.Range("M105").Value = -788
$ws.Range("N105").Value = -8339660.5
$ws.Range("H134").Value = 3086.5945
$ws.Range("I134").Value = 3205.9714
$ws.Range("J134").Value = 997.5
$ws.Range("K134").Value = 9617.914199999999
$ws.Range("L134").Value = 2992.5
$ws.Range("M134").Value = -7082.914199999999
$ws.Range("N134").Value = -8062.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2555.4092
$ws.Range("I31").Value = 1046.2
$ws.Range("K31").Value = 1046.2
$ws.Range("M31").Value = -751.2
$ws.Range("H34").Value = 2555.4092
$ws.Range("I34").Value = 1046.2
$ws.Range("K34").Value = 1046.2
$ws.Range("M34").Value = -844.2
$ws.Range("H58").Value = 18282.467
$ws.Range("I58").Value = 1482.375
$ws.Range("J58").Value = 37482.57
$ws.Range("K58").Value = 1482.375
$ws.Range("L58").Value = 37482.57
$ws.Range("M58").Value = -1279.375
$ws.Range("N58").Value = -37888.57
$ws.Range("H99").Value = 5816.6665
$ws.Range("I99").Value = 3633.3333
$ws.Range("J99").Value = 8000
$ws.Range("K99").Value = 3633.3333
$ws.Range("L99").Value = 8000
$ws.Range("M99").Value = -2135.3333
$ws.Range("N99").Value = -10996
$ws.Range("H126").Value = 5816.6665
$ws.Range("I126").Value = 3633.3333
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 10899.9999
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -8429.999899999999
$ws.Range("N126").Value = -28940
$ws.Range("H132").Value = 4800
$ws.Range("I132").Value = 2081
$ws.Range("J132").Value = 6975.2
$ws.Range("K132").Value = 6243
$ws.Range("L132").Value = 20925.6
$ws.Range("M132").Value = -3713
$ws.Range("N132").Value = -25985.6
$ws.Range("H134").Value = 1229
$ws.Range("I134").Value = 1160.3889
$ws.Range("J134").Value = 1434.8334
$ws.Range("K134").Value = 3481.1667
$ws.Range("L134").Value = 4304.5002
$ws.Range("M134").Value = -946.1666999999998
$ws.Range("N134").Value = -9374.5002
$ws.Range("H136").Value = 18282.467
$ws.Range("I136").Value = 1482.375
$ws.Range("J136").Value = 37482.57
$ws.Range("K136").Value = 4447.125
$ws.Range("L136").Value = 112447.71
$ws.Range("M136").Value = -1897.125
$ws.Range("N136").Value = -117547.71

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 300
$ws.Range("I44").Value = 300
$ws.Range("K44").Value = 900
$ws.Range("M44").Value = -502
$ws.Range("H45").Value = 539.1667
$ws.Range("I45").Value = 308.75
$ws.Range("K45").Value = 926.25
$ws.Range("M45").Value = -394.25
$ws.Range("H122").Value = 486.5263
$ws.Range("J122").Value = 548.53845
$ws.Range("L122").Value = 4936.84605
$ws.Range("N122").Value = -9836.84605
$ws.Range("H131").Value = 710.23
$ws.Range("J131").Value = 721.61053
$ws.Range("L131").Value = 2164.83159
$ws.Range("N131").Value = -12244.83159

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 17911.064
$ws.Range("I132").Value = 1249.1666
$ws.Range("J132").Value = 40981.383
$ws.Range("K132").Value = 3747.4998
$ws.Range("L132").Value = 122944.149
$ws.Range("M132").Value = -1217.4998
$ws.Range("N132").Value = -128004.149

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 838.2
$ws.Range("I16").Value = 747.75
$ws.Range("K16").Value = 747.75
$ws.Range("M16").Value = -577.75
$ws.Range("H22").Value = 8333.333000000001
$ws.Range("I22").Value = 10001
$ws.Range("J22").Value = 7499.5
$ws.Range("K22").Value = 10001
$ws.Range("L22").Value = 7499.5
$ws.Range("M22").Value = -9706
$ws.Range("N22").Value = -8089.5
$ws.Range("H27").Value = 8333.333000000001
$ws.Range("I27").Value = 10001
$ws.Range("J27").Value = 7499.5
$ws.Range("K27").Value = 10001
$ws.Range("L27").Value = 7499.5
$ws.Range("M27").Value = -9894
$ws.Range("N27").Value = -7713.5
$ws.Range("H34").Value = 70024
$ws.Range("J34").Value = 70024
$ws.Range("L34").Value = 70024
$ws.Range("N34").Value = -70368
$ws.Range("H100").Value = 2159.8
$ws.Range("I100").Value = 1750
$ws.Range("J100").Value = 2433
$ws.Range("K100").Value = 1750
$ws.Range("L100").Value = 2433
$ws.Range("M100").Value = -1209
$ws.Range("N100").Value = -3515
$ws.Range("H122").Value = 1511434.9
$ws.Range("I122").Value = 1963485.2
$ws.Range("J122").Value = 4600
$ws.Range("K122").Value = 5890455.6
$ws.Range("L122").Value = 13800
$ws.Range("M122").Value = -5888005.6
$ws.Range("N122").Value = -18700
$ws.Range("H132").Value = 2636.1177
$ws.Range("I132").Value = 2216.6155
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 6649.8465
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -4119.8465
$ws.Range("N132").Value = -17058.5
$ws.Range("H136").Value = 2109.6316
$ws.Range("I136").Value = 1967.9231
$ws.Range("J136").Value = 2416.6667
$ws.Range("K136").Value = 5903.7693
$ws.Range("L136").Value = 7250.000100000001
$ws.Range("M136").Value = -3353.7693
$ws.Range("N136").Value = -12350.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1150.375
$ws.Range("I122").Value = 1115.8462
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 3347.5386
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -897.5385999999999
$ws.Range("N122").Value = -8800
$ws.Range("H132").Value = 1166.8214
$ws.Range("I132").Value = 680.05
$ws.Range("K132").Value = 2040.15
$ws.Range("M132").Value = 489.8500000000001
